$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 227; this shifts the existing rows 227-277 down to 228-278
# and extends the used range / dimension to A1:R278 automatically.
$ws.Rows.Item(227).Insert()

# Populate the newly inserted row 227 with the new weekly record.
$ws.Cells.Item(227, 1).Value = 8
$ws.Cells.Item(227, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(227, 3).Value = "Coquimbo"
$ws.Cells.Item(227, 4).Value = 45135
$ws.Cells.Item(227, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(227, 5).Value = 4
$ws.Cells.Item(227, 6).Value = 100112001
$ws.Cells.Item(227, 7).Value = "Berenjena"
$ws.Cells.Item(227, 8).Value = "Sin especificar"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 440
$ws.Cells.Item(227, 11).Value = 8000
$ws.Cells.Item(227, 12).Value = 9000
$ws.Cells.Item(227, 13).Value = 8500
$ws.Cells.Item(227, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(227, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(227, 16).Value = 170
$ws.Cells.Item(227, 17).Value = 50
$ws.Cells.Item(227, 18).Value = "Hortaliza"
